$d = $word.ActiveDocument

$replacements = @(
    @{old="56×54="; new="37×51="},
    @{old="52×16="; new="17×44="},
    @{old="15×77="; new="25×13="},
    @{old="69×42="; new="35×21="},
    @{old="96×42="; new="71×59="},
    @{old="12×23="; new="89×49="},
    @{old="31×20="; new="92×24="},
    @{old="61×81="; new="19×45="},
    @{old="25×49="; new="42×86="},
    @{old="16×82="; new="19×53="},
    @{old="33×30="; new="42×88="},
    @{old="11×51="; new="65×71="},
    @{old="47×73="; new="15×75="},
    @{old="57×55="; new="66×22="},
    @{old="28×56="; new="83×38="},
    @{old="60×12="; new="68×74="},
    @{old="59×12="; new="65×96="},
    @{old="42×41="; new="36×64="},
    @{old="83×20="; new="72×28="},
    @{old="70×70="; new="86×11="},
    @{old="84×11="; new="76×76="},
    @{old="46×93="; new="17×31="},
    @{old="55×62="; new="28×51="},
    @{old="98×86="; new="88×57="},
    @{old="96×90="; new="16×73="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
